$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Actualización" footnote text (B38's "Nota" text is unchanged)
$ws.Range("B37").Value = "Actualización: Mayo 2025."

# Updated 2022 (N) and 2023 (O) GDP-by-state figures
$ws.Range("N5").Value = 314602.24699999997
$ws.Range("O5").Value = 334294.52299999999
$ws.Range("N6").Value = 919300.27800000005
$ws.Range("O6").Value = 951029.42799999996
$ws.Range("N7").Value = 178129.82399999999
$ws.Range("O7").Value = 181120.095
$ws.Range("N8").Value = 472680.02399999998
$ws.Range("O8").Value = 514860.22200000001
$ws.Range("N9").Value = 935841.223
$ws.Range("O9").Value = 929917.34499999997
$ws.Range("N10").Value = 148438.43400000001
$ws.Range("O10").Value = 154588.88
$ws.Range("N11").Value = 374518.00300000003
$ws.Range("O11").Value = 379736.48300000001
$ws.Range("N12").Value = 911271.14199999999
$ws.Range("O12").Value = 941368.24300000002
$ws.Range("N13").Value = 3650377.6189999999
$ws.Range("O13").Value = 3806101.5189999999
$ws.Range("N14").Value = 300478.11300000001
$ws.Range("O14").Value = 310240.19099999999
$ws.Range("N15").Value = 1110918.71
$ws.Range("O15").Value = 1128468.281
$ws.Range("N16").Value = 310105.89899999998
$ws.Range("O16").Value = 313074.75
$ws.Range("N17").Value = 431251.88199999998
$ws.Range("O17").Value = 441422.255
$ws.Range("N18").Value = 1797547.0630000001
$ws.Range("O18").Value = 1849411.733
$ws.Range("N19").Value = 2218452.727
$ws.Range("O19").Value = 2269428.8089999999
$ws.Range("N20").Value = 652066.83900000004
$ws.Range("O20").Value = 678519.16099999996
$ws.Range("N21").Value = 260119.46299999999
$ws.Range("O21").Value = 265778.52899999998
$ws.Range("N22").Value = 162601.04500000001
$ws.Range("O22").Value = 162761.44
$ws.Range("N23").Value = 1928658.429
$ws.Range("O23").Value = 1995895.4169999999
$ws.Range("N24").Value = 415502.71100000001
$ws.Range("O24").Value = 446836.29499999998
$ws.Range("N25").Value = 825802.78200000001
$ws.Range("O25").Value = 851623.35699999996
$ws.Range("N26").Value = 589221.23699999996
$ws.Range("O26").Value = 613730.68500000006
$ws.Range("N27").Value = 356727.57299999997
$ws.Range("O27").Value = 406186.435
$ws.Range("N28").Value = 532793.15
$ws.Range("O28").Value = 575114.47400000005
$ws.Range("N29").Value = 516319.26899999997
$ws.Range("O29").Value = 517548.272
$ws.Range("N30").Value = 790043.70400000003
$ws.Range("O30").Value = 829550.01599999995
$ws.Range("N31").Value = 637775.22699999996
$ws.Range("O31").Value = 653095.64099999995
$ws.Range("N32").Value = 748872.57
$ws.Range("O32").Value = 749222.91899999999
$ws.Range("N33").Value = 146641.59899999999
$ws.Range("O33").Value = 148529.16200000001
$ws.Range("N34").Value = 1041978.086
$ws.Range("O34").Value = 1069284.605
$ws.Range("N35").Value = 369043.26799999998
$ws.Range("O35").Value = 389502.13699999999
$ws.Range("N36").Value = 225013.34899999999
$ws.Range("O36").Value = 228860.43
